$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Introduction ")
$c = $ws1.Cells.Item(7,4)
$c.Value = 1.7
$v2 = $c.Value2
Write-Output $v2
